$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.138.20'
$ws.Cells.Item(2, 5).Value = '  +2.21%  '
$ws.Cells.Item(3, 4).Value = '1.809.84'
$ws.Cells.Item(3, 5).Value = '  +3.03%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  -1.25%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '329.59'
$ws.Cells.Item(5, 5).Value = '  +1.49%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9990'
$ws.Cells.Item(6, 5).Value = '  -1.09%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4436'
$ws.Cells.Item(7, 5).Value = '  +4.29%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3724'
$ws.Cells.Item(8, 5).Value = '  +3.24%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '44.70'
$ws.Cells.Item(9, 5).Value = '  -0.83%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.07705'
$ws.Cells.Item(10, 5).Value = '  +4.75%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.122'
$ws.Cells.Item(11, 5).Value = '  +0.49%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.0000'
$ws.Cells.Item(12, 5).Value = '  -1.26%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '21.98'
$ws.Cells.Item(13, 5).Value = '  +1.76%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.301'
$ws.Cells.Item(14, 5).Value = '  +3.16%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.474'
$ws.Cells.Item(15, 5).Value = '  +3.31%  '
$ws.Cells.Item(16, 4).Value = '1.814.70'
$ws.Cells.Item(16, 5).Value = '  +2.90%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '93.56'
$ws.Cells.Item(17, 5).Value = '  +12.63%  '
$ws.Cells.Item(18, 5).Value = '  +2.01%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06470'
$ws.Cells.Item(19, 5).Value = '  +7.63%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '1.001'
$ws.Cells.Item(20, 5).Value = '  -0.79%  '
$ws.Cells.Item(21, 5).Value = '  +3.85%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.258'
$ws.Cells.Item(22, 5).Value = '  +2.65%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.5340'
$ws.Cells.Item(23, 5).Value = '  -2.00%  '
$ws.Cells.Item(24, 4).Value = '28.190.19'
$ws.Cells.Item(24, 5).Value = '  +2.22%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '11.71'
$ws.Cells.Item(25, 5).Value = '  +4.02%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.139'
$ws.Cells.Item(26, 5).Value = '  -10.97%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.60'
$ws.Cells.Item(27, 5).Value = '  +3.41%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '155.58'
$ws.Cells.Item(28, 5).Value = '  +3.02%  '
$ws.Cells.Item(29, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.330'
$ws.Cells.Item(29, 5).Value = '  -1.73%  '
$ws.Cells.Item(30, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(30, 4).Value = '2.017.88'
$ws.Cells.Item(30, 5).Value = '  +2.53%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '127.58'
$ws.Cells.Item(31, 5).Value = '  +0.50%  '
$ws.Cells.Item(32, 5).Value = '  -5.79%  '
$ws.Cells.Item(33, 5).Value = '  +5.42%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.09221'
$ws.Cells.Item(34, 5).Value = '  +1.84%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.663'
$ws.Cells.Item(35, 5).Value = '  -0.70%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '13.07'
$ws.Cells.Item(36, 5).Value = '  +4.74%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02341'
$ws.Cells.Item(37, 5).Value = '  +3.42%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2170'
$ws.Cells.Item(38, 5).Value = '  +0.45%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.168'
$ws.Cells.Item(39, 5).Value = '  +3.19%  '
$ws.Cells.Item(40, 2).Value = 'TheSandbox'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.6564'
$ws.Cells.Item(40, 5).Value = '  +1.88%  '
$ws.Cells.Item(41, 2).Value = 'Hedera'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.06195'
$ws.Cells.Item(41, 5).Value = '  +0.96%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.192'
$ws.Cells.Item(42, 5).Value = '  +1.71%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '8.066'
$ws.Cells.Item(43, 5).Value = '  +1.96%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.9997'
$ws.Cells.Item(44, 5).Value = '  -1.08%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.91'
$ws.Cells.Item(45, 5).Value = '  +0.40%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.388'
$ws.Cells.Item(46, 5).Value = '  -2.90%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.6069'
$ws.Cells.Item(47, 5).Value = '  +3.04%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.767'
$ws.Cells.Item(48, 5).Value = '  -0.32%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '126.70'
$ws.Cells.Item(49, 5).Value = '  +1.34%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.030'
$ws.Cells.Item(50, 5).Value = '  +4.43%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06977'
$ws.Cells.Item(51, 5).Value = '  +1.42%  '
